$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency price/volume figures (and restore the
# correct row order for two coin pairs that got swapped) to match the
# latest run of the GitHub Actions scraper.
#
# Price/volume values are plain text in the workbook (e.g. "30.339.50",
# "1.000"), so for any new value that Excel would otherwise auto-convert
# into a number (losing formatting such as trailing zeros) we force the
# cell to Text format first, preserving the exact literal string.

$ws.Range("D2").Value = "30.339.50"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "1.939.99"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.74"
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7194"
$ws.Range("E6").Value = "  -7.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3347"
$ws.Range("E8").Value = "  -5.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.65"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07274"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8154"
$ws.Range("E11").Value = "  -5.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08154"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "1.938.66"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.542"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.16"
$ws.Range("E15").Value = "  -5.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.86"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("D17").Value = "30.350.18"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "254.20"
$ws.Range("E18").Value = "  -7.34%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008247"
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.885"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "2.194.75"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.869"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.36"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.430"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.47"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1325"
$ws.Range("E29").Value = "  -10.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.559"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.348"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.466"
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.236"
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05227"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.272"
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7550"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.737"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02004"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.846"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.42"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.673"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4578"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.029"
$ws.Range("E43").Value = "  -5.64%  "
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.69"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.862"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.414"
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.89"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4168"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.501"
$ws.Range("E51").Value = "  -0.23%  "
